# "Generate Report for Handback" -- populate the "Latest Target File" /
# "Latest Handback File" / "Latest Handback DateTime" columns for the two
# localized-language sheets (zh-cn, de-de) now that handback has completed,
# flip the Overview/per-language "Status" column from "Ready for handoff" to
# "Handed back: in sync with en-US", and widen the columns that now hold the
# longer file-name / status text so it isn't clipped.

$wb = $excel.ActiveWorkbook

# Cornflower-blue (#6495ED) underlined font used by the existing hyperlink
# cells (A2/A3) -- OLE/VBA color value is BGR-packed: 0xED*65536 + 0x95*256 + 0x64
$hyperlinkColor = 15570276

$statusText = "Handed back: in sync with en-US"

$rows = @(
    @{ Row = 2; Guid = "a70c4802-e0a0-44ae-9b0f-23258a963346" },
    @{ Row = 3; Guid = "d5427d45-234e-468a-975d-98b7678503d4" }
)

# Per-language sheet info: the xlf hash segment embedded in the handoff/handback
# file names, and the handback timestamp recorded for that language.
$languages = @(
    @{ Sheet = "zh-cn"; Hash = "97a1bb77022f50eae0ac4b8dd31c4374efde7b51"; Hash2 = "a1146ec33dbc7dc5cfa7ac92583271269efad4ce"; HandbackTime = "2016-08-30 09:30:54" },
    @{ Sheet = "de-de"; Hash = "97a1bb77022f50eae0ac4b8dd31c4374efde7b51"; Hash2 = "a1146ec33dbc7dc5cfa7ac92583271269efad4ce"; HandbackTime = "2016-08-30 09:31:08" }
)

foreach ($lang in $languages) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    foreach ($r in $rows) {
        $row = $r.Row
        $guid = $r.Guid

        # Status column (C): handoff -> handed back
        $ws.Range("C$row").Value = $statusText

        # Latest Target File (I): the source .md file, same display text/target
        # as the Source File Name hyperlink in column A
        $mdName = "$guid.md"
        $ws.Range("I$row").Value = $mdName

        # Latest Handback File (J): the translated xlf that came back
        $xlfHash = if ($row -eq 2) { $lang.Hash } else { $lang.Hash2 }
        $ws.Range("J$row").Value = "$guid.$xlfHash.$($lang.Sheet).xlf"

        # Latest Handback DateTime (K): when the handback happened
        $ws.Range("K$row").Value = $lang.HandbackTime

        # Give I$row the same hyperlink (and hyperlink styling) as A$row already has
        $targetUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/10cd4fcdcd07f8f8aa39b1545a4acaef2bb17982/e2e/$mdName"
        $ws.Hyperlinks.Add($ws.Range("I$row"), $targetUrl, "", "", $mdName)
        $ws.Range("I$row").Font.Underline = $true
        $ws.Range("I$row").Font.Color = $hyperlinkColor
    }

    # Widen Status (C) and the two newly-populated columns (I, J) so the
    # longer handback text/file names aren't clipped.
    $ws.Columns.Item(3).ColumnWidth = 29 + 1/6   # -> stored width 29.9777047293527 (closest representable)
    $ws.Columns.Item(9).ColumnWidth = 40 - 5/6   # -> stored width 40
    $ws.Columns.Item(10).ColumnWidth = 40 - 5/6  # -> stored width 40
}

# Overview sheet mirrors the per-language Status text in its zh-cn/de-de
# columns (E/F) and needs the matching width bump.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $statusText
$overview.Range("F2").Value = $statusText
$overview.Range("E3").Value = $statusText
$overview.Range("F3").Value = $statusText
$overview.Columns.Item(5).ColumnWidth = 29 + 1/6
$overview.Columns.Item(6).ColumnWidth = 29 + 1/6

Write-Host "Handback report generated"
